$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the component name and manufacturer for row 2 (test data for the
# database, per the "add pandas and openpyxl" test entries).
$ws.Range("A2").Value = "Teste banco de dados"
$ws.Range("C2").Value = "teste banco de dados"

# Move the active selection to C19, matching the author's last cursor
# position when the workbook was saved.
$ws.Range("C19").Select()
